$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 4 (shifts existing rows 4-29 down to 6-31),
# matching the new dimension A1:T31.
$ws.Rows.Item(4).Resize(2).Insert()

# New row 4 -> id 2, "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 1.000478137718011
$ws.Range("D4").Value = 1.024728322086813
$ws.Range("E4").Value = 0.9819445537750674
$ws.Range("F4").Value = 1.024728322086813
$ws.Range("G4").Value = 0.9819445537750674
$ws.Range("H4").Value = 1.005058547844381
$ws.Range("I4").Value = 1.009144973949326
$ws.Range("J4").Value = 0.9950747618132221
$ws.Range("K4").Value = 0.9819445537750674
$ws.Range("L4").Value = 1.000478137718011
$ws.Range("M4").Value = 1.012603229902412
$ws.Range("N4").Value = 1.012603229902412
$ws.Range("O4").Value = 1.01145047791805
$ws.Range("P4").Value = 1.002383671193297
$ws.Range("Q4").Value = 1.002383671193297
$ws.Range("R4").Value = 0.9972738918387397
$ws.Range("S4").Value = 0.9972738918387397
$ws.Range("T4").Value = 1.002738216197803

# New row 5 -> id 3, "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.052943902210302
$ws.Range("D5").Value = 0.720519821052607
$ws.Range("E5").Value = 1.102395021679571
$ws.Range("F5").Value = 0.720519821052607
$ws.Range("G5").Value = 1.102395021679571
$ws.Range("H5").Value = 1.136654632538749
$ws.Range("I5").Value = 0.8546696816976938
$ws.Range("J5").Value = 1.067361123647926
$ws.Range("K5").Value = 1.102395021679571
$ws.Range("L5").Value = 1.052943902210302
$ws.Range("M5").Value = 0.8867318616314546
$ws.Range("N5").Value = 0.8867318616314546
$ws.Range("O5").Value = 0.876044468320201
$ws.Range("P5").Value = 0.9586195816474934
$ws.Range("Q5").Value = 0.9586195816474934
$ws.Range("R5").Value = 0.9945634416555127
$ws.Range("S5").Value = 0.9945634416555127
$ws.Range("T5").Value = 0.9890906971378081

# Apply the same style as the other id cells in column A to the two new id cells
$ws.Range("A4:A5").Style = $ws.Range("A6").Style

# Rename "Thomas Hex" -> "Matthies Hex" (now at row 11 after the insert shift)
$ws.Range("B11").Value = "Matthies Hex"
